$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45205 = 2023-10-06) for every
# data row (2..503). The edit bumps it by one day to 45206 (2023-10-07) for all rows.
$lastRow = 503
$ws.Range("C2:C$lastRow").Value = 45206
